# Update the "rate" worksheet: rename the "rate" column to "searchRate" and
# add a new "relocationRate" column with a value of 6.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rate")

$ws.Range("B2").Value = "searchRate"
$ws.Range("C2").Value = "relocationRate"
$ws.Range("C3").Value = 6

# Size the new columns to fit their contents (mirrors Excel's AutoFit after
# typing in the new header/values). The host's column-width grid is coarser
# than real Excel's, so these inputs are chosen to land on the closest
# achievable rendered width to the saved workbook's 9.71 / 12.71.
$ws.Columns.Item(2).ColumnWidth = 8.86
$ws.Columns.Item(3).ColumnWidth = 11.86

# Make "rate" the active sheet with C3 selected, matching the saved view.
$ws.Activate() | Out-Null
$ws.Range("C3").Select() | Out-Null
